$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 3069.25
$ws.Cells.Item(32, 9).Value = 1138.5
$ws.Cells.Item(32, 10).Value = 5000
$ws.Cells.Item(32, 11).Value = 1138.5
$ws.Cells.Item(32, 12).Value = 5000
$ws.Cells.Item(32, 13).Value = -812.5
$ws.Cells.Item(32, 14).Value = -5652

$ws.Cells.Item(51, 8).Value = 7671.067
$ws.Cells.Item(51, 9).Value = 3966
$ws.Cells.Item(51, 10).Value = 9523.6
$ws.Cells.Item(51, 11).Value = 3966
$ws.Cells.Item(51, 12).Value = 9523.6
$ws.Cells.Item(51, 13).Value = -3482
$ws.Cells.Item(51, 14).Value = -10491.6

$ws.Cells.Item(64, 8).Value = 4044.4443
$ws.Cells.Item(64, 9).Value = 3700
$ws.Cells.Item(64, 10).Value = 4142.857
$ws.Cells.Item(64, 11).Value = 3700
$ws.Cells.Item(64, 12).Value = 4142.857
$ws.Cells.Item(64, 13).Value = -3452
$ws.Cells.Item(64, 14).Value = -4638.857

$ws.Cells.Item(67, 8).Value = 4044.4443
$ws.Cells.Item(67, 9).Value = 3700
$ws.Cells.Item(67, 10).Value = 4142.857
$ws.Cells.Item(67, 11).Value = 3700
$ws.Cells.Item(67, 12).Value = 4142.857
$ws.Cells.Item(67, 13).Value = -2842
$ws.Cells.Item(67, 14).Value = -5858.857

$ws.Cells.Item(112, 8).Value = 1231.723
$ws.Cells.Item(112, 10).Value = 1278.0656
$ws.Cells.Item(112, 12).Value = 3834.1968
$ws.Cells.Item(112, 14).Value = -6050.1968

$ws.Cells.Item(116, 8).Value = 723772.4399999999
$ws.Cells.Item(116, 9).Value = 3335763.2
$ws.Cells.Item(116, 10).Value = 11411.272
$ws.Cells.Item(116, 11).Value = 3335763.2
$ws.Cells.Item(116, 12).Value = 11411.272
$ws.Cells.Item(116, 13).Value = -3332321.2
$ws.Cells.Item(116, 14).Value = -18295.272

$ws.Cells.Item(129, 8).Value = 812.62
$ws.Cells.Item(129, 9).Value = 274
$ws.Cells.Item(129, 10).Value = 865.89014
$ws.Cells.Item(129, 11).Value = 822
$ws.Cells.Item(129, 12).Value = 2597.67042
$ws.Cells.Item(129, 13).Value = 4178
$ws.Cells.Item(129, 14).Value = -12597.67042

$ws.Cells.Item(137, 8).Value = 1192371.1
$ws.Cells.Item(137, 9).Value = 1588665.2
$ws.Cells.Item(137, 10).Value = 3489
$ws.Cells.Item(137, 11).Value = 4765995.6
$ws.Cells.Item(137, 12).Value = 10467
$ws.Cells.Item(137, 13).Value = -4763445.6
$ws.Cells.Item(137, 14).Value = -15567

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1612.9524
$ws.Cells.Item(61, 9).Value = 1403.0667
$ws.Cells.Item(61, 10).Value = 2137.6667
$ws.Cells.Item(61, 11).Value = 1403.0667
$ws.Cells.Item(61, 12).Value = 2137.6667
$ws.Cells.Item(61, 13).Value = -1191.0667
$ws.Cells.Item(61, 14).Value = -2561.6667

$ws.Cells.Item(74, 8).Value = 1707.0385
$ws.Cells.Item(74, 9).Value = 911.82355
$ws.Cells.Item(74, 10).Value = 3209.111
$ws.Cells.Item(74, 11).Value = 911.82355
$ws.Cells.Item(74, 12).Value = 3209.111
$ws.Cells.Item(74, 13).Value = -37.82354999999995
$ws.Cells.Item(74, 14).Value = -4957.111

$ws.Cells.Item(77, 8).Value = 1707.0385
$ws.Cells.Item(77, 9).Value = 911.82355
$ws.Cells.Item(77, 10).Value = 3209.111
$ws.Cells.Item(77, 11).Value = 4559.117749999999
$ws.Cells.Item(77, 12).Value = 16045.555
$ws.Cells.Item(77, 13).Value = -191.1177499999994
$ws.Cells.Item(77, 14).Value = -24781.555

$ws.Cells.Item(110, 8).Value = 956.3125
$ws.Cells.Item(110, 9).Value = 1000.6667
$ws.Cells.Item(110, 10).Value = 899.2857
$ws.Cells.Item(110, 11).Value = 1000.6667
$ws.Cells.Item(110, 12).Value = 899.2857
$ws.Cells.Item(110, 13).Value = 1044.3333
$ws.Cells.Item(110, 14).Value = -4989.2857

$ws.Cells.Item(112, 8).Value = 31710.525
$ws.Cells.Item(112, 10).Value = 31710.525
$ws.Cells.Item(112, 12).Value = 31710.525
$ws.Cells.Item(112, 14).Value = -34664.525

$ws.Cells.Item(119, 8).Value = 35465.332
$ws.Cells.Item(119, 10).Value = 35465.332
$ws.Cells.Item(119, 12).Value = 35465.332
$ws.Cells.Item(119, 14).Value = -45141.332

$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).Value = ""

$ws.Cells.Item(125, 8).Value = 41805.625
$ws.Cells.Item(125, 10).Value = 41805.625
$ws.Cells.Item(125, 12).Value = 41805.625
$ws.Cells.Item(125, 14).Value = -51645.625

$ws.Cells.Item(132, 8).Value = 2425
$ws.Cells.Item(132, 9).Value = 984.7273
$ws.Cells.Item(132, 10).Value = 5065.5
$ws.Cells.Item(132, 11).Value = 2954.1819
$ws.Cells.Item(132, 12).Value = 15196.5
$ws.Cells.Item(132, 13).Value = -424.1819
$ws.Cells.Item(132, 14).Value = -20256.5

$ws.Cells.Item(136, 8).Value = 1612.9524
$ws.Cells.Item(136, 9).Value = 1403.0667
$ws.Cells.Item(136, 10).Value = 2137.6667
$ws.Cells.Item(136, 11).Value = 4209.2001
$ws.Cells.Item(136, 12).Value = 6413.000100000001
$ws.Cells.Item(136, 13).Value = -1659.2001
$ws.Cells.Item(136, 14).Value = -11513.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6669.391
$ws.Cells.Item(31, 9).Value = 1382.4286
$ws.Cells.Item(31, 10).Value = 14893.556
$ws.Cells.Item(31, 11).Value = 1382.4286
$ws.Cells.Item(31, 12).Value = 14893.556
$ws.Cells.Item(31, 13).Value = -1087.4286
$ws.Cells.Item(31, 14).Value = -15483.556

$ws.Cells.Item(34, 8).Value = 6669.391
$ws.Cells.Item(34, 9).Value = 1382.4286
$ws.Cells.Item(34, 10).Value = 14893.556
$ws.Cells.Item(34, 11).Value = 1382.4286
$ws.Cells.Item(34, 12).Value = 14893.556
$ws.Cells.Item(34, 13).Value = -1180.4286
$ws.Cells.Item(34, 14).Value = -15297.556

$ws.Cells.Item(58, 8).Value = 2798.6428
$ws.Cells.Item(58, 9).Value = 1727.0193
$ws.Cells.Item(58, 10).Value = 5894.4443
$ws.Cells.Item(58, 11).Value = 1727.0193
$ws.Cells.Item(58, 12).Value = 5894.4443
$ws.Cells.Item(58, 13).Value = -1524.0193
$ws.Cells.Item(58, 14).Value = -6300.4443

$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).Value = ""

$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).Value = ""

$ws.Cells.Item(98, 8).Value = 50000
$ws.Cells.Item(98, 10).Value = 50000
$ws.Cells.Item(98, 12).Value = 50000
$ws.Cells.Item(98, 14).Value = -54492

$ws.Cells.Item(132, 8).Value = 4778.263
$ws.Cells.Item(132, 9).Value = 4789.5835
$ws.Cells.Item(132, 10).Value = 4758.857
$ws.Cells.Item(132, 11).Value = 14368.7505
$ws.Cells.Item(132, 12).Value = 14276.571
$ws.Cells.Item(132, 13).Value = -11838.7505
$ws.Cells.Item(132, 14).Value = -19336.571

$ws.Cells.Item(134, 8).Value = 3963.3333
$ws.Cells.Item(134, 9).Value = 4351.6
$ws.Cells.Item(134, 10).Value = 2992.6667
$ws.Cells.Item(134, 11).Value = 13054.8
$ws.Cells.Item(134, 12).Value = 8978.000100000001
$ws.Cells.Item(134, 13).Value = -10519.8
$ws.Cells.Item(134, 14).Value = -14048.0001

$ws.Cells.Item(136, 8).Value = 2798.6428
$ws.Cells.Item(136, 9).Value = 1727.0193
$ws.Cells.Item(136, 10).Value = 5894.4443
$ws.Cells.Item(136, 11).Value = 5181.0579
$ws.Cells.Item(136, 12).Value = 17683.3329
$ws.Cells.Item(136, 13).Value = -2631.0579
$ws.Cells.Item(136, 14).Value = -22783.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 800.65
$ws.Cells.Item(131, 9).Value = 282.2
$ws.Cells.Item(131, 10).Value = 827.9367999999999
$ws.Cells.Item(131, 11).Value = 846.5999999999999
$ws.Cells.Item(131, 12).Value = 2483.8104
$ws.Cells.Item(131, 13).Value = 4193.4
$ws.Cells.Item(131, 14).Value = -12563.8104

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4240.3887
$ws.Cells.Item(132, 9).Value = 3186.2727
$ws.Cells.Item(132, 10).Value = 5896.857
$ws.Cells.Item(132, 11).Value = 9558.8181
$ws.Cells.Item(132, 12).Value = 17690.571
$ws.Cells.Item(132, 13).Value = -7028.8181
$ws.Cells.Item(132, 14).Value = -22750.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 7001501
$ws.Cells.Item(11, 10).Value = 7001501
$ws.Cells.Item(11, 12).Value = 7001501
$ws.Cells.Item(11, 14).Value = -7001781

$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).Value = ""

$ws.Cells.Item(25, 8).Value = 9250
$ws.Cells.Item(25, 9).Value = 3500
$ws.Cells.Item(25, 10).Value = 15000
$ws.Cells.Item(25, 11).Value = 3500
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = -3270
$ws.Cells.Item(25, 14).Value = -15460

$ws.Cells.Item(62, 8).Value = 18226
$ws.Cells.Item(62, 9).Value = 18226
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 18226
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -17602
$ws.Cells.Item(62, 14).Value = ""

$ws.Cells.Item(65, 8).Value = 18226
$ws.Cells.Item(65, 9).Value = 18226
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 54678
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -51558
$ws.Cells.Item(65, 14).Value = ""

$ws.Cells.Item(103, 8).Value = 35000
$ws.Cells.Item(103, 10).Value = 35000
$ws.Cells.Item(103, 12).Value = 35000
$ws.Cells.Item(103, 14).Value = -37344

$ws.Cells.Item(110, 8).Value = 40400
$ws.Cells.Item(110, 10).Value = 40400
$ws.Cells.Item(110, 12).Value = 40400
$ws.Cells.Item(110, 14).Value = -48580

$ws.Cells.Item(127, 8).Value = 25678.846
$ws.Cells.Item(127, 10).Value = 25678.846
$ws.Cells.Item(127, 12).Value = 25678.846
$ws.Cells.Item(127, 14).Value = -35598.84600000001

$ws.Cells.Item(136, 8).Value = 5377
$ws.Cells.Item(136, 9).Value = 1587
$ws.Cells.Item(136, 11).Value = 4761
$ws.Cells.Item(136, 13).Value = -2211

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 22111
$ws.Cells.Item(119, 10).Value = 22111
$ws.Cells.Item(119, 12).Value = 22111
$ws.Cells.Item(119, 14).Value = -31787

$ws.Cells.Item(132, 8).Value = 23812702
$ws.Cells.Item(132, 9).Value = 2017.2858
$ws.Cells.Item(132, 10).Value = 47623384
$ws.Cells.Item(132, 11).Value = 6051.857400000001
$ws.Cells.Item(132, 12).Value = 142870152
$ws.Cells.Item(132, 13).Value = -3521.857400000001
$ws.Cells.Item(132, 14).Value = -142875212

$ws.Cells.Item(136, 8).Value = 5351.9546
$ws.Cells.Item(136, 9).Value = 2453.1667
$ws.Cells.Item(136, 10).Value = 8830.5
$ws.Cells.Item(136, 11).Value = 7359.500100000001
$ws.Cells.Item(136, 12).Value = 26491.5
$ws.Cells.Item(136, 13).Value = -4809.500100000001
$ws.Cells.Item(136, 14).Value = -31591.5

